# Commit: "added funder to grants, removed some that didn't get funded"
# Change the "type" column (A) from "rejected" to "omit" for the grants
# that were not actually funded / should be excluded from the CV listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows (1-indexed, matching worksheet rows) whose type changes from
# "rejected" to "omit":
#   15 - Reproducible collaboration: Interfacing between Microsoft Word and R Markdown
#   17 - Putting Large-Scale Data to Work in Applied Educational Settings
#   18 - Open and Reproducible Research in Education
#   19 - County Moderators of District Expenditures on Student Achievement
$rowsToOmit = @(15, 17, 18, 19)

foreach ($r in $rowsToOmit) {
    $ws.Cells.Item($r, 1).Value = "omit"
}

# Update the selected cell to reflect where the editor left off.
$ws.Range("A20").Select()
